# Atualização automática de preços de eletricidade
# Updates row 2 (the single data row) of the SpotPTTable with the latest
# hourly spot-price values and the recomputed summary-slot columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46023

$ws.Range("B2").Value = 108.88
$ws.Range("C2").Value = 103.6
$ws.Range("D2").Value = 99.29000000000001
$ws.Range("E2").Value = 88.47
$ws.Range("F2").Value = 78.93000000000001
$ws.Range("G2").Value = 71.42
$ws.Range("H2").Value = 71.22
$ws.Range("I2").Value = 75
$ws.Range("J2").Value = 67.92
$ws.Range("K2").Value = 61.59
$ws.Range("L2").Value = 56.61
$ws.Range("M2").Value = 52.8
$ws.Range("N2").Value = 54.26
$ws.Range("O2").Value = 57.58
$ws.Range("P2").Value = 61.4
$ws.Range("Q2").Value = 64.92
$ws.Range("R2").Value = 77.29000000000001
$ws.Range("S2").Value = 93.76000000000001
$ws.Range("T2").Value = 101.18
$ws.Range("U2").Value = 104.71
$ws.Range("V2").Value = 106.78
$ws.Range("W2").Value = 106.35
$ws.Range("X2").Value = 104.16
$ws.Range("Y2").Value = 96.68000000000001
$ws.Range("Z2").Value = 81.87

$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 103.49
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 106.56
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 106.24
$ws.Range("AG2").Value = "4h-16h"
